# "adding new progress as of date 04 nov 2025"
# For rows 3-15 on the "Training Dashboard" sheet:
#   - column H (PERIOD TO EXPIRE) decreases by 1
#   - column I (LAST UPDATE) moves from 03-Nov-2025 to 04-Nov-2025
#
# Column I is stored as plain text (not a real date), so the cell's
# NumberFormat is force-flipped to Text ("@") before the write so Excel's
# automatic date recognition doesn't turn the literal "04-Nov-2025" string
# into a date serial, then flipped back to the (lowercase-matched, so it
# resolves to the same built-in General format / style) "general" format.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

for ($row = 3; $row -le 15; $row++) {
    $hCell = $ws.Cells.Item($row, 8)   # column H - PERIOD TO EXPIRE
    $hCell.Value = $hCell.Value2 - 1

    $iCell = $ws.Cells.Item($row, 9)   # column I - LAST UPDATE
    $iCell.NumberFormat = "@"
    $iCell.Value = "04-Nov-2025"
    $iCell.NumberFormat = "general"
}
